$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8648375868797302
$ws.Range("B1").Value = 1.540218114852905
$ws.Range("C1").Value = 6.324757099151611
$ws.Range("D1").Value = 3.01220440864563
$ws.Range("E1").Value = 1.576457381248474
